$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$generatedCells = @("A3","A49","A95","A131","A166","A189","A215","A224","A274")
foreach ($c in $generatedCells) {
    $ws.Range($c).Value = "Generated: 05/22/2024, 03:59 PM"
}

$pricingDateCells = @("A10","A56","A102","A138","A172","A195","A221","A231","A280")
foreach ($c in $pricingDateCells) {
    $ws.Range($c).Value = "Pricing Date: 05/22/2024"
}

# COST ADJUSTMENT - OTHER total increased 64882 -> 79421, rippling through
# dependent cached totals elsewhere in the report (+14539 each).
$ws.Range("G88").Value = 79421
$ws.Range("G90").Value = 335762.48

$ws.Range("B174").Value = 408042.46
$ws.Range("B176").Value = 408042.46
$ws.Range("C183").Value = 408042.46
$ws.Range("C184").Value = 408042.46

$ws.Range("C263").Value = 93137
$ws.Range("B265").Value = 79421
$ws.Range("C269").Value = 93158.76
$ws.Range("B272").Value = 93158.76
